$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new shared strings in the same order they appear in the target
# workbook (Hydraulic Diameter, Throat, Chamber, Exit, Average, Other Analyses)
# so that the shared-strings table indices line up.
$ws.Range("R2").Value = "Hydraulic Diameter"
$ws.Range("R4").Value = "Throat"
$ws.Range("R3").Value = "Chamber"
$ws.Range("R5").Value = "Exit"
$ws.Range("R6").Value = "Average"
$ws.Range("R1").Value = "Other Analyses"

$ws.Range("S3").Formula = "=4*1000^2*I6/(2*M3+2*M6)"
$ws.Range("T3").Value = "mm"

$ws.Range("S4").Formula = "=4*1000^2*I6/(2*M4+2*M7)"
$ws.Range("T4").Value = "mm"

$ws.Range("S5").Formula = "=4*1000^2*I6/(2*M5+2*M8)"
$ws.Range("T5").Value = "mm"

$ws.Range("S6").Formula = "=AVERAGE(S3:S5)"

$ws.Range("S3:S6").NumberFormat = "0.00000"

$ws.Range("R2").Select()
